# The "Fix sorting issues (use depth stencil gubbins)" task has been
# resolved (depth stencil added to fix z-sorting), so remove its row
# from the ToDo list. This shifts all subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 currently holds the "Fix sorting issues..." task (A2/B2).
# Select the whole row first (mirrors a user right-clicking the row
# header and choosing Delete), then delete it so the rows below shift up.
$ws.Rows(2).Select()
$ws.Rows(2).Delete()
